$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.06436660403446501, 0.02161267820907678, -0.02475332161006712),
    @(-0.002483299958681726, -0.07227465348399185, -0.05542938006312945),
    @(0.05864306564604772, -0.01641368349690153, -0.03955352419744367),
    @(0.0103050321340562, 0.04435412786172135, 0.04177123000440386),
    @(0.02360463045213534, 0.01676559415848359, 0.06694286187057909),
    @(0.02699094848788307, -0.01045774781833528, 0.02604145086977801),
    @(-0.01472052470173514, 0.001075654007170485, -0.02027806955511146),
    @(-0.0279470856260994, -0.001367806902398176, -0.02752877472211478),
    @(-0.02199114867202608, 0.004966600231178349, -0.03471972080676452),
    @(-0.01893682202891156, 0.0004515091119253129, 0.08639095043358573),
    @(-0.007516298443078927, 0.01816660165786737, 0.06252737019373024),
    @(-0.01159979990157086, 0.006480483672298237, 0.04596096941310422),
    @(-0.01335271759687552, 0.001062374369686679, 0.02341207367894439),
    @(-0.02131388465995362, -0.005225553865665965, 0.007270624132259934),
    @(-0.03625352340547929, -0.01064366328975433, 0.008558753068032388),
    @(-0.02540402729874076, 0.01100221381563203, 0.005378270172513961),
    @(0.0394008085131645, -0.0215329993516206, -0.001527163083665),
    @(-0.005351710416700385, -0.006526962327568362, -0.02914889580468922),
    @(0.0001261571584187555, 0.01954104799939229, 0.05001127266365517),
    @(0.1110513180737254, 0.04867666449559752, -0.1655312067140687),
    @(0.1914996167887811, 0.03030422819859344, 0.02057685541069637),
    @(0.1882859338884768, 0.0655750582480559, 0.4253946024438608),
    @(0.06060181622919789, 0.2589603908683945, 0.4051430983387903),
    @(-0.2129994086597282, 0.3793273995751927, 0.01590242139671223),
    @(0.5540348317312142, -0.6639242237028846, -1.121634874006974),
    @(-0.05442012457743672, -0.6851783446643664, -1.266947746276855),
    @(-1.690450101442949, 0.1374446801517266, -2.446953586910074),
    @(-0.7898089613603607, 0.3619974125986522, -1.304489494017943),
    @(0.258422552243524, 0.7682891954546356, -0.810219791272407),
    @(0.04960623909921762, 0.7135171605193199, -0.9448558195777682)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 2, 1).Value = $row[0]
    $ws.Cells.Item($i + 2, 2).Value = $row[1]
    $ws.Cells.Item($i + 2, 3).Value = $row[2]
}

$ws.Range("A1:C31").Select()